# Added surefireplugin and firefox browser support
#
# On the DATA sheet:
#   - insert a new "browser" column (all rows set to "chrome") between the
#     existing "execute" and "username" columns
#   - fix row 4's "execute" flag, which should be "no" (it was mistakenly
#     left as "yes")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Insert a new column before column C; this shifts the old C:E -> D:F
$ws.Range("C1").EntireColumn.Insert()

# Populate the new "browser" column
$ws.Range("C1").Value = "browser"
$ws.Range("C2:C6").Value = "chrome"

# Correct the "execute" value for the test1/loginLogoutTest row
$ws.Range("B4").Value = "no"

# Leave the selection where the author left it when they saved the file
[void]$ws.Range("F8").Select()
